# Update the System Vision, User Requirements and Software Requirements
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# --- Project title (B1) ---
$ws.Range("B1").Value = "Nutrient Analyzer"

# --- Activity names (column B, rows 5-39) ---
# Written in explicit row order (5 -> 39), except row 32 is poked in just
# ahead of row 31 so new shared-string entries get interned in the same
# order the source workbook used ("8.2 ..." before "8.1 ...").
$ws.Cells.Item(5, 2).Value = "1. Project Objectives Documentation"
$ws.Cells.Item(6, 2).Value = "2. Project Scope Definition"
$ws.Cells.Item(7, 2).Value = "3. Stakeholder Identification"
$ws.Cells.Item(8, 2).Value = "4.Design"
$ws.Cells.Item(9, 2).Value = "4.1 Design Documentation"
$ws.Cells.Item(10, 2).Value = "4.2 User Interface (UI) Design"
$ws.Cells.Item(11, 2).Value = "4.3 User Experience (UX) Design"
$ws.Cells.Item(12, 2).Value = "5. Data Preprocessing"
$ws.Cells.Item(13, 2).Value = "5.1Nutritional_Food_Database.csv Import"
$ws.Cells.Item(14, 2).Value = "5.2 Data Cleaning and Validation"
$ws.Cells.Item(15, 2).Value = "5.3 Data Integration"
$ws.Cells.Item(16, 2).Value = "6. Development"
$ws.Cells.Item(17, 2).Value = "6.1 Desktop Application Development"
$ws.Cells.Item(18, 2).Value = "6.2 Feature Implementation"
$ws.Cells.Item(19, 2).Value = "6.2.1 Food Search"
$ws.Cells.Item(20, 2).Value = "6.2.2 Nutrition Breakdown"
$ws.Cells.Item(21, 2).Value = "6.2.3 Nutrition Range Filter"
$ws.Cells.Item(22, 2).Value = "6.2.4 Nutrition Level Filter"
$ws.Cells.Item(23, 2).Value = "6.2.5 Additional Feature"
$ws.Cells.Item(24, 2).Value = "6.2.6 Graphical User Interface (GUI) Development"
$ws.Cells.Item(25, 2).Value = "7 Testing and Validation"
$ws.Cells.Item(26, 2).Value = "7.1 Unit Testing"
$ws.Cells.Item(27, 2).Value = "7.2 Integration Testing"
$ws.Cells.Item(28, 2).Value = "7.3 User Acceptance Testing (UAT)"
$ws.Cells.Item(29, 2).Value = "7.4 Usability Testing"
$ws.Cells.Item(30, 2).Value = "8. Documentation"
$ws.Cells.Item(32, 2).Value = "8.2 Design Documentation"
$ws.Cells.Item(31, 2).Value = "8.1 Project Plan Documentation"
$ws.Cells.Item(33, 2).Value = "8.3 User Manual"
$ws.Cells.Item(34, 2).Value = "9. Project Closure"
$ws.Cells.Item(35, 2).Value = "9.1 Final Report"
$ws.Cells.Item(36, 2).Value = "9.2 Lessons Learned"
$ws.Cells.Item(37, 2).Value = "9.3 Project Evaluation"
$ws.Cells.Item(38, 2).Value = "Activity 34"
$ws.Cells.Item(39, 2).Value = "Activity 35"

# --- Plan start / duration values for the first couple of activities ---
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = 0

$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 0

# --- View state: zoom + selection ---
$ws.Activate()
$ws.Range("C7").Select()
$excel.ActiveWindow.Zoom = 85
